$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear the "proceed" (column C) values ---
# Rows that still show an (empty) cell afterwards - i.e. just the cell
# contents were cleared, keeping the cell/style in place.
$ws.Range("C2:C8").ClearContents()
$ws.Range("C14").ClearContents()

# Rows where the cell disappears completely (full clear, not just contents).
$ws.Range("C9:C13").Clear()
$ws.Range("C15:C23").Clear()

# --- Add "Molecule ChEMBL ID" to J18 / J20 ---
$ws.Range("J18").Value = "Molecule ChEMBL ID"
$ws.Range("J20").Value = "Molecule ChEMBL ID"

# --- Append new row 24 (BIOCODES / BIO_AbaumanniiMIC_ChEMBL dataset) ---
$ws.Range("A24").Value = "BIOCODES"
$ws.Range("E24").Value = "BIO_ChemBL_AbaumanniiMIC_firstprocessing.csv"
$ws.Range("B24").Value = "BIO_AbaumanniiMIC_ChEMBL"
$ws.Range("T24").Value = "BIO_AbaumanniiMIC_ChEMBL"
$ws.Range("C24").Value = "yes"
$ws.Range("F24").Value = "curation8"
$ws.Range("H24").Value = "curation4"
$ws.Range("I24").Value = "ID"
$ws.Range("J24").Value = "Molecule ChEMBL ID"
$ws.Range("Q24").Value = "y"
$ws.Range("R24").Value = "nM"

$ws.Range("B24").Copy()
$ws.Range("T24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H32").Select()
